$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E3")

$data = @(
    @(2, 57.7, 0.88, 2.54, 26),
    @(3, 61.5, 0.66, 3, 26),
    @(4, 67.90000000000001, 0.97, 2.12, 28),
    @(5, 66.7, 1.06, 3.29, 27),
    @(6, 88, 0.88, 4.11, 25),
    @(7, 80, 1.08, 3.11, 15),
    @(8, 16.7, 0.64, 4.62, 6),
    @(9, 88.5, 1.28, 1.92, 26),
    @(10, 65.40000000000001, 1.49, 3.01, 26),
    @(11, 71.40000000000001, 1.77, 1.74, 28),
    @(12, 33.3, 2.29, 3.05, 12),
    @(13, 87.5, 2.6, 1.6, 24),
    @(14, 70.40000000000001, 1.87, 2.15, 27),
    @(15, 56, 2.32, 4.4, 25),
    @(16, 73.3, 3.44, 12.8, 15),
    @(17, 100, 5.45, 13.11, 6),
    @(18, 58.8, 1.9, 1.72, 17),
    @(19, 50, 1.61, 4.19, 24),
    @(20, 84.2, 1.36, 2.32, 19),
    @(21, 100, 1.56, 2.92, 24),
    @(22, 40, 2.48, 12.72, 15),
    @(23, 56.2, 1.3, 9.74, 16),
    @(24, 57.9, 2.24, 11.05, 19),
    @(25, 100, 2.79, 11.21, 8),
    @(26, 100, 1.38, 8.970000000000001, 9),
    @(27, 100, 1.23, 6.73, 9),
    @(28, 90.90000000000001, 1.17, 4.59, 11),
    @(29, 81.5, 1.25, 7.35, 27),
    @(30, 89.3, 2.68, 6.59, 28),
    @(31, 84, 0.99, 7.89, 25),
    @(32, 92, 2.05, 4.53, 25),
    @(33, 100, 0.91, 3.67, 14),
    @(34, 100, 1.19, 4.33, 24),
    @(35, 100, 1.46, 5.68, 21),
    @(36, 100, 1.27, 4.9, 18),
    @(37, 100, 1.43, 13.82, 7),
    @(38, 10.5, 0, 0.76, 38),
    @(39, 71.40000000000001, 0, 0.08, 14),
    @(40, 100, 0, 0.67, 8)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 8).Value = $row[4]
}
